$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per the diff. Numeric-looking text values are
# prefixed with a literal leading apostrophe so Excel keeps them as
# text (matching the original inline-string cell content) instead of
# silently converting them to numbers.

$ws.Range('D2').Value = '57.480.36'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '3.109.33'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''525.23'
$ws.Range('E5').Value = '  +0.44%  '
$ws.Range('D6').Value = '''136.89'
$ws.Range('E6').Value = '  -2.78%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.107.84'
$ws.Range('E8').Value = '  +0.30%  '
$ws.Range('E9').Value = '  +2.36%  '
$ws.Range('D10').Value = '''7.24'
$ws.Range('E10').Value = '  +0.64%  '
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('E12').Value = '  +3.28%  '
$ws.Range('D13').Value = '3.645.24'
$ws.Range('D15').Value = '''25.32'
$ws.Range('E15').Value = '  -2.65%  '
$ws.Range('D16').Value = '''0.0000163'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('D17').Value = '57.617.25'
$ws.Range('E17').Value = '  +0.42%  '
$ws.Range('D18').Value = '3.107.13'
$ws.Range('E18').Value = '  +0.26%  '
$ws.Range('E19').Value = '  -2.60%  '
$ws.Range('D20').Value = '''12.45'
$ws.Range('E20').Value = '  -2.74%  '
$ws.Range('D21').Value = '''7.90'
$ws.Range('E21').Value = '  -1.87%  '
$ws.Range('D22').Value = '''347.80'
$ws.Range('E22').Value = '  +3.03%  '
$ws.Range('D23').Value = '''5.81'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').Value = '''68.15'
$ws.Range('E25').Value = '  +2.23%  '
$ws.Range('E26').Value = '  -1.91%  '
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('D30').Value = '''7.42'
$ws.Range('E30').Value = '  +3.53%  '
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('E32').Value = '  +0.22%  '
$ws.Range('E33').Value = '  -7.54%  '
$ws.Range('D34').Value = '''20.97'
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('D35').Value = '''4.97'
$ws.Range('E35').Value = '  +7.52%  '
$ws.Range('E36').Value = '  -2.12%  '
$ws.Range('D37').Value = '''158.08'
$ws.Range('E37').Value = '  +0.77%  '
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('D39').Value = '''25.87'
$ws.Range('E39').Value = '  -4.72%  '
$ws.Range('D40').Value = '''1.23'
$ws.Range('E40').Value = '  -3.30%  '
$ws.Range('E41').Value = '  +6.84%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').Value = '''0.0665'
$ws.Range('E42').Value = '  +0.94%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '''1.61'
$ws.Range('E43').Value = '  +6.56%  '
$ws.Range('D44').Value = '''0.700'
$ws.Range('E44').Value = '  +2.25%  '
$ws.Range('D45').Value = '3.148.05'
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('D46').Value = '2.356.26'
$ws.Range('E46').Value = '  +1.99%  '
$ws.Range('D47').Value = '''36.48'
$ws.Range('E47').Value = '  -0.26%  '
$ws.Range('D48').Value = '''1.00'
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('D49').Value = '''0.0267'
$ws.Range('E49').Value = '  +3.34%  '
$ws.Range('D50').Value = '''0.959'
$ws.Range('E50').Value = '  -1.14%  '
$ws.Range('E51').Value = '  +0.55%  '
